$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 487 F/G values (AUC fix)
$ws.Cells.Item(487, 6).Value = 45924.50271662037
$ws.Cells.Item(487, 7).Value = 45924.50266359954

# Append new experiment rows 488-527
$rows = @(
    @(488, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.62051621528, 45924.6204903125),
    @(489, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.65678487268, 45924.65678231481),
    @(490, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.6660155787, 45924.66601310185),
    @(491, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.6660155787, 45924.66601310185),
    @(492, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.66860956018, 45924.66860641204),
    @(493, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.66860956018, 45924.66860641204),
    @(494, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.67399969907, 45924.67399782407),
    @(495, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.67574260417, 45924.67574130787),
    @(496, 'Fucntionality_test_BZR_with_SVC_Simple-Prototype-GED_poly', 'BZR', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_BZR.joblib', 45924.67574260417, 45924.67574130787),
    @(497, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45924.67655055555, 45924.67654515046),
    @(498, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45924.68070023148, 45924.68069513889),
    @(499, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.51793018518, 45925.51783008102),
    @(500, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.52300373842, 45925.52189148148),
    @(501, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.52649340278, 45925.52648737268),
    @(502, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.53015548611, 45925.53014549769),
    @(503, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.54304729166, 45925.54304013889),
    @(504, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.54304729166, 45925.54304013889),
    @(505, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.54536806713, 45925.54536305556),
    @(506, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.54850828704, 45925.54850324074),
    @(507, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.55185165509, 45925.55184623843),
    @(508, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45925.55185165509, 45925.55184623843),
    @(509, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.57046309028, 45925.57028105324),
    @(510, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.57171554398, 45925.57147372685),
    @(511, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.57743530093, 45925.57673060185),
    @(512, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.58201896991, 45925.58174378472),
    @(513, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.58601043982, 45925.58583722222),
    @(514, 'Fucntionality_test_BZR_with_SVC_Random-Walk-Edit_precomputed', 'BZR', 'SVC_Random-Walk-Edit_precomputed', 'SVC_Random-Walk-Edit_precomputed_trained_on_BZR.joblib', 45925.58601043982, 45925.58583722222),
    @(515, 'Fucntionality_test_Letter-high_with_(7)-NN_Classifier_GED', 'Letter-high', '(7)-NN_Classifier_GED', '(7)-NN_Classifier_GED_trained_on_Letter-high.joblib', 45926.55947094908, 45926.55942363426),
    @(516, 'Fucntionality_test_MUTAG_with_(7)-NN_Classifier_GED', 'MUTAG', '(7)-NN_Classifier_GED', '(7)-NN_Classifier_GED_trained_on_MUTAG.joblib', 45926.56548694445, 45926.56548652778),
    @(517, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45926.56598181713, 45926.56597805556),
    @(518, 'Fucntionality_test_Letter-high_with_(7)-NN_Classifier_GED', 'Letter-high', '(7)-NN_Classifier_GED', '(7)-NN_Classifier_GED_trained_on_Letter-high.joblib', 45926.57018256944, 45926.57013594908),
    @(519, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.60692478009, 45926.60692107639),
    @(520, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.60778174768, 45926.60777612268),
    @(521, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.60778174768, 45926.60777612268),
    @(522, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.63602863426, 45926.63602075232),
    @(523, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.6371171412, 45926.63710895833),
    @(524, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.65641547454, 45926.65641295139),
    @(525, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.66002355324, 45926.66002107639),
    @(526, 'Fucntionality_test_BZR_with_SVC_Zero-GED_precomputed', 'BZR', 'SVC_Zero-GED_precomputed', 'SVC_Zero-GED_precomputed_trained_on_BZR.joblib', 45926.66002355324, 45926.66002107639),
    @(527, 'Fucntionality_test_MUTAG_with_SVC_Trivial-GED_precomputed', 'MUTAG', 'SVC_Trivial-GED_precomputed', 'SVC_Trivial-GED_precomputed_trained_on_MUTAG.joblib', 45926.82500684206, 45926.82500614467)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
